$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "image_0_2.png"
$ws.Range("C2").Value = 21
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = "3 inch white drizzle bone"
$ws.Range("G2").Value = "ABC123"

$ws.Range("C3").Value = 15
